$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text updated: "Ready for handoff" -> "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# zh-cn row: refreshed handback datetime, and error detail cleared (version is now up to date)
$wsZhCn.Range("K2").Value = "2016-09-02 22:55:20"
$wsZhCn.Range("P2").Value = ""

# de-de row: refreshed handback datetime, and error detail cleared (version is now up to date)
$wsDeDe.Range("K2").Value = "2016-09-02 22:55:27"
$wsDeDe.Range("P2").Value = ""

# Column width adjustments (widened to better fit longer status text)
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
